$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row labels: "_old" -> "_FV2404", "_new" -> "_FV2410"
for ($col = 1; $col -le 10; $col++) {
  $cell = $ws.Cells.Item(1, $col)
  $v = $cell.Value()
  $cell.Value = $v -replace "_old$", "_FV2404"
}
for ($col = 12; $col -le 21; $col++) {
  $cell = $ws.Cells.Item(1, $col)
  $v = $cell.Value()
  $cell.Value = $v -replace "_new$", "_FV2410"
}

# Turn the used range into a native Excel Table ("Table1")
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U68"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (row 1)
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
Write-Host "Applied header renames, table, and frozen header row."
